# RPA datasets push 2024-04-28
# Remove the "하나31호스팩" (Hana 31st SPAC) record, which is the last
# data row on each of the three sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Sheet "01_리그테이블": last data row is row 12 (A12:Q12)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(12).Delete()

# Sheet "02_통합집계_Rawdata": last data row is row 11 (A11:T11)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(11).Delete()

# Sheet "03_IPO현황_Summary": the 하나31호스팩 record sits on row 11
# (A11:L11); deleting it shifts the following row (삼현) up to row 11.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(11).Delete()
